$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Total" header in column T (row 1)
$ws.Range("T1").Value = "Total"

# Add Total column values for existing rows (2-6)
$ws.Range("T2").Value = 88810
$ws.Range("T3").Value = 10472
$ws.Range("T4").Value = 41529
$ws.Range("T5").Value = 13943
$ws.Range("T6").Value = 53907

# Add new row 7: "Outros"
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 6818
$ws.Range("C7").Value = 287
$ws.Range("D7").Value = 451
$ws.Range("E7").Value = 2073
$ws.Range("F7").Value = 2742
$ws.Range("G7").Value = 2825
$ws.Range("H7").Value = 3232
$ws.Range("I7").Value = 3730
$ws.Range("J7").Value = 3811
$ws.Range("K7").Value = 4486
$ws.Range("L7").Value = 5309
$ws.Range("M7").Value = 5692
$ws.Range("N7").Value = 5710
$ws.Range("O7").Value = 5809
$ws.Range("P7").Value = 5764
$ws.Range("Q7").Value = 6791
$ws.Range("R7").Value = 21645
$ws.Range("S7").Value = 523
$ws.Range("T7").Value = 87698

# Add new row 8: "Total"
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 7681
$ws.Range("C8").Value = 485
$ws.Range("D8").Value = 680
$ws.Range("E8").Value = 2585
$ws.Range("F8").Value = 3526
$ws.Range("G8").Value = 3958
$ws.Range("H8").Value = 5052
$ws.Range("I8").Value = 6589
$ws.Range("J8").Value = 8167
$ws.Range("K8").Value = 11516
$ws.Range("L8").Value = 16327
$ws.Range("M8").Value = 21536
$ws.Range("N8").Value = 25829
$ws.Range("O8").Value = 28804
$ws.Range("P8").Value = 29554
$ws.Range("Q8").Value = 32949
$ws.Range("R8").Value = 90412
$ws.Range("S8").Value = 709
$ws.Range("T8").Value = 296359
